# Auto-generated edit script: updates cryptos price (D) and volume/1h (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.272.80"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").Value = "2.608.68"
$ws.Range("E3").Value = "  +3.89%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.21"
$ws.Range("E5").Value = "  +3.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.38"
$ws.Range("E6").Value = "  +3.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.604"
$ws.Range("E7").Value = "  +2.61%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.581"
$ws.Range("E9").Value = "  +6.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.48"
$ws.Range("E10").Value = "  +7.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").Value = "  +5.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.26"
$ws.Range("E13").Value = "  +7.00%  "
$ws.Range("D14").Value = "3.010.91"
$ws.Range("E14").Value = "  +4.29%  "
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("D16").Value = "2.602.30"
$ws.Range("E16").Value = "  +3.83%  "
$ws.Range("E17").Value = "  +5.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.92"
$ws.Range("E18").Value = "  +2.89%  "
$ws.Range("D19").Value = "46.471.12"
$ws.Range("E19").Value = "  +1.94%  "
$ws.Range("E20").Value = "  +4.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.94"
$ws.Range("E21").Value = "  -2.42%  "
$ws.Range("E22").Value = "  +4.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.48"
$ws.Range("E23").Value = "  +4.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "272.65"
$ws.Range("E24").Value = "  +9.62%  "
$ws.Range("E25").Value = "  +6.22%  "
$ws.Range("E26").Value = "  +6.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.14"
$ws.Range("E27").Value = "  +26.62%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.03"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.60"
$ws.Range("E30").Value = "  +5.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "38.98"
$ws.Range("E31").Value = "  -2.95%  "
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("E33").Value = "  +10.02%  "
$ws.Range("E35").Value = "  +2.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.23"
$ws.Range("E36").Value = "  +4.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0839"
$ws.Range("E37").Value = "  +4.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "151.07"
$ws.Range("E38").Value = "  +1.44%  "
$ws.Range("E39").Value = "  +5.12%  "
$ws.Range("E40").Value = "  +4.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.39"
$ws.Range("E41").Value = "  +40.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.89"
$ws.Range("E42").Value = "  +2.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.65"
$ws.Range("E43").Value = "  +8.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0332"
$ws.Range("E44").Value = "  +6.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.09"
$ws.Range("E45").Value = "  +0.71%  "
$ws.Range("D46").Value = "2.126.75"
$ws.Range("E46").Value = "  +6.24%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "93.39"
$ws.Range("E48").Value = "  +2.84%  "
$ws.Range("E49").Value = "  +8.48%  "
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "109.19"
$ws.Range("E51").Value = "  +3.71%  "
